$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.839.80'
$ws.Range('E2').Value = '  +2.47%  '
$ws.Range('D3').Value = '2.565.44'
$ws.Range('E3').Value = '  +1.94%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.39'
$ws.Range('E5').Value = '  +1.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.73'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('D9').Value = '2.564.51'
$ws.Range('E9').Value = '  +1.90%  '
$ws.Range('E10').Value = '  +10.78%  '
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('E14').Value = '  +6.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.28'
$ws.Range('E15').Value = '  +1.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000181'
$ws.Range('E16').Value = '  +4.89%  '
$ws.Range('D17').Value = '69.742.23'
$ws.Range('E17').Value = '  +2.59%  '
$ws.Range('D18').Value = '2.571.80'
$ws.Range('E18').Value = '  +3.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.69'
$ws.Range('E19').Value = '  +1.85%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '366.59'
$ws.Range('E20').Value = '  +3.96%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.19'
$ws.Range('E21').Value = '  +1.30%  '
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.75'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.31'
$ws.Range('E25').Value = '  -0.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.78'
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  +1.19%  '
$ws.Range('D28').Value = '2.690.20'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('D30').Value = '0.0₃0925'
$ws.Range('E30').Value = '  +0.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '519.84'
$ws.Range('E31').Value = '  +2.03%  '
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.28'
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('E34').Value = '  +1.47%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('E36').Value = '  -1.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '163.07'
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.02'
$ws.Range('E38').Value = '  +3.18%  '
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('E41').Value = '  +1.00%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.96'
$ws.Range('E43').Value = '  +1.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.327'
$ws.Range('E44').Value = '  -1.18%  '
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.05'
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '153.54'
$ws.Range('E47').Value = '  +4.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.64'
$ws.Range('E48').Value = '  +2.23%  '
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('D50').Value = '0.0₆0258'
$ws.Range('E50').Value = '  -0.36%  '
$ws.Range('E51').Value = '  +1.76%  '
